$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace B6's URL (growing-up stage-4) with the anti-reflux URL that was in B7
$antiReflux = $ws.Cells.Item(7, 2).Value2
$ws.Range("B6").Value = $antiReflux

# Delete row 7 entirely (shifts nothing below it up since it's the last row)
$ws.Rows("7").Delete()
